$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G2: change numeric 1 to text "1;2"
$ws.Range("G2").Value = "1;2"

# D3: set to "POSL"
$ws.Range("D3").Value = "POSL"

# E4: set to "OPEN"
$ws.Range("E4").Value = "OPEN"

# Update selection to G3 as in the final sheetView
$ws.Range("G3").Select()
